# Trade #121 closed at 2026-02-17 09:29:07 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the new closed trade (row 122) to
# both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Helper: write a literal text value into a cell without letting
# Excel's smart "looks like a date" auto-conversion turn strings such
# as "2026-02-17" into a date serial number. Leading apostrophe forces
# text entry; ClearFormats() afterwards drops the quote-prefix /
# number-format style that the apostrophe trick leaves behind so the
# cell keeps the workbook's default (unstyled) look.
# ---------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.79               # Current Capital
$summary.Range("B4").Value = 0.8                   # Total P&L $
$summary.Range("B5").Value = 0.13                  # Total P&L %
$summary.Range("B6").Value = 121                    # Total Trades
$summary.Range("B7").Value = 55                    # Winning Trades
$summary.Range("B9").Value = 45.45                 # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.79                 # Capital
$status.Range("D4").Value = 121                     # Trades
$status.Range("E4").Value = 0.8                    # P&L $
$status.Range("F4").Value = 0.79                   # P&L %
$status.Range("G4").Value = 45.45                  # Win Rate %

# ---------------------------------------------------------------
# New closed trade (row 122) shared by "All Trades" and "MarketMaking"
# ---------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 122

    $ws.Cells.Item($row, 1).Value = 121                 # A: Trade #
    Set-TextValue $ws.Cells.Item($row, 2) "2026-02-17"  # B: Date
    Set-TextValue $ws.Cells.Item($row, 3) "09:29:00"    # C: Time
    Set-TextValue $ws.Cells.Item($row, 4) "MarketMaking" # D: Strategy
    Set-TextValue $ws.Cells.Item($row, 5) "DOWN"        # E: Side
    $ws.Cells.Item($row, 6).Value = 0.83                # F: Entry Price
    $ws.Cells.Item($row, 7).Value = 0.9399999999999999  # G: Exit Price
    Set-TextValue $ws.Cells.Item($row, 8) "CLOSED"      # H: Status
    $ws.Cells.Item($row, 9).Value = 13.253              # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0.11               # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100.79             # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                  # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                  # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                # N: Confidence
    Set-TextValue $ws.Cells.Item($row, 15) "Normal spread capture: 19600 bps" # O: Entry Reason
    Set-TextValue $ws.Cells.Item($row, 16) "early_exit" # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.15               # Q: Duration (min)
}
